$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 148.875
$ws.Range("I9").Value = 145.28572
$ws.Range("J9").Value = 174
$ws.Range("K9").Value = 145.28572
$ws.Range("L9").Value = 174
$ws.Range("M9").Value = 23.71428
$ws.Range("N9").Value = -512

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2478
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 2597.5
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 2597.5
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2947.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3105.52
$ws.Range("I86").Value = 2710.2856
$ws.Range("J86").Value = 3608.5454
$ws.Range("K86").Value = 2710.2856
$ws.Range("L86").Value = 3608.5454
$ws.Range("M86").Value = -1587.2856
$ws.Range("N86").Value = -5854.5454

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3105.52
$ws.Range("I89").Value = 2710.2856
$ws.Range("J89").Value = 3608.5454
$ws.Range("K89").Value = 13551.428
$ws.Range("L89").Value = 18042.727
$ws.Range("M89").Value = -7935.428
$ws.Range("N89").Value = -29274.727

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3377.3333
$ws.Range("I113").Value = 1200
$ws.Range("J113").Value = 3649.5
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 3649.5
$ws.Range("M113").Value = 2054
$ws.Range("N113").Value = -10157.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1905
$ws.Range("J121").Value = 1905
$ws.Range("L121").Value = 5715
$ws.Range("N121").Value = -9209

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1330.68
$ws.Range("J138").Value = 1698.125
$ws.Range("L138").Value = 5094.375
$ws.Range("N138").Value = -15374.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4318.3
$ws.Range("I32").Value = 4412.4644
$ws.Range("K32").Value = 4412.4644
$ws.Range("M32").Value = -4125.4644

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1622.5333
$ws.Range("I122").Value = 1660.7858
$ws.Range("J122").Value = 1087
$ws.Range("K122").Value = 4982.357400000001
$ws.Range("L122").Value = 3261
$ws.Range("M122").Value = -2532.357400000001
$ws.Range("N122").Value = -8161

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4862.2
$ws.Range("I132").Value = 4828
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 14484
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -11954
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 16932.3
$ws.Range("J135").Value = 16932.3
$ws.Range("L135").Value = 16932.3
$ws.Range("N135").Value = -27072.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3184.9333
$ws.Range("I86").Value = 3616.7144
$ws.Range("J86").Value = 2177.4443
$ws.Range("K86").Value = 3616.7144
$ws.Range("L86").Value = 2177.4443
$ws.Range("M86").Value = -2493.7144
$ws.Range("N86").Value = -4423.4443

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3184.9333
$ws.Range("I89").Value = 3616.7144
$ws.Range("J89").Value = 2177.4443
$ws.Range("K89").Value = 18083.572
$ws.Range("L89").Value = 10887.2215
$ws.Range("M89").Value = -12467.572
$ws.Range("N89").Value = -22119.2215

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 142859900
$ws.Range("I105").Value = 200002860
$ws.Range("K105").Value = 200002860
$ws.Range("M105").Value = -200001113

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 10970
$ws.Range("I134").Value = 1833.875
$ws.Range("J134").Value = 35333
$ws.Range("K134").Value = 5501.625
$ws.Range("L134").Value = 105999
$ws.Range("M134").Value = -2966.625
$ws.Range("N134").Value = -111069

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2286.16
$ws.Range("I31").Value = 1225.8
$ws.Range("K31").Value = 1225.8
$ws.Range("M31").Value = -930.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2286.16
$ws.Range("I34").Value = 1225.8
$ws.Range("K34").Value = 1225.8
$ws.Range("M34").Value = -1023.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 398.63635
$ws.Range("I105").Value = 383.125
$ws.Range("J105").Value = 440
$ws.Range("K105").Value = 383.125
$ws.Range("L105").Value = 440
$ws.Range("M105").Value = 1363.875
$ws.Range("N105").Value = -3934

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2746.1
$ws.Range("I132").Value = 2208.8572
$ws.Range("J132").Value = 3999.6667
$ws.Range("K132").Value = 6626.571599999999
$ws.Range("L132").Value = 11999.0001
$ws.Range("M132").Value = -4096.571599999999
$ws.Range("N132").Value = -17059.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 17544864
$ws.Range("I134").Value = 20834332
$ws.Range("J134").Value = 1041.3334
$ws.Range("K134").Value = 62502996
$ws.Range("L134").Value = 3124.0002
$ws.Range("M134").Value = -62500461
$ws.Range("N134").Value = -8194.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2075.7
$ws.Range("J69").Value = 2245.2222
$ws.Range("L69").Value = 6735.6666
$ws.Range("N69").Value = -8357.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 2075.7
$ws.Range("J72").Value = 2245.2222
$ws.Range("L72").Value = 20206.9998
$ws.Range("N72").Value = -28318.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 749.6875
$ws.Range("J113").Value = 749.6875
$ws.Range("L113").Value = 2249.0625
$ws.Range("N113").Value = -6589.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 408.3
$ws.Range("I121").Value = 342.55554
$ws.Range("J121").Value = 1000
$ws.Range("K121").Value = 1027.66662
$ws.Range("L121").Value = 3000
$ws.Range("M121").Value = 282.33338
$ws.Range("N121").Value = -5620

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 29413312
$ws.Range("I131").Value = 200000240
$ws.Range("J131").Value = 1771.4482
$ws.Range("K131").Value = 600000720
$ws.Range("L131").Value = 5314.3446
$ws.Range("M131").Value = -599995680
$ws.Range("N131").Value = -15394.3446

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3003.6667
$ws.Range("I132").Value = 2006
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 6018
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -3488
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1366.8214
$ws.Range("I22").Value = 1382.1111
$ws.Range("J22").Value = 1339.3
$ws.Range("K22").Value = 1382.1111
$ws.Range("L22").Value = 1339.3
$ws.Range("M22").Value = -1087.1111
$ws.Range("N22").Value = -1929.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1366.8214
$ws.Range("I27").Value = 1382.1111
$ws.Range("J27").Value = 1339.3
$ws.Range("K27").Value = 1382.1111
$ws.Range("L27").Value = 1339.3
$ws.Range("M27").Value = -1275.1111
$ws.Range("N27").Value = -1553.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 40000
$ws.Range("J36").Value = 40000
$ws.Range("L36").Value = 40000
$ws.Range("N36").Value = -41124

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5591
$ws.Range("I46").Value = 2400.2
$ws.Range("J46").Value = 8250
$ws.Range("K46").Value = 2400.2
$ws.Range("L46").Value = 8250
$ws.Range("M46").Value = -2212.2
$ws.Range("N46").Value = -8626

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 324.625
$ws.Range("I55").Value = 85.28570999999999
$ws.Range("K55").Value = 85.28570999999999
$ws.Range("M55").Value = 87.71429000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1707.875
$ws.Range("I61").Value = 1734.8
$ws.Range("J61").Value = 1663
$ws.Range("K61").Value = 1734.8
$ws.Range("L61").Value = 1663
$ws.Range("M61").Value = -1532.8
$ws.Range("N61").Value = -2067

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2131.3125
$ws.Range("I100").Value = 1854.6364
$ws.Range("J100").Value = 2740
$ws.Range("K100").Value = 1854.6364
$ws.Range("L100").Value = 2740
$ws.Range("M100").Value = -1313.6364
$ws.Range("N100").Value = -3822

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1707.875
$ws.Range("I113").Value = 1734.8
$ws.Range("J113").Value = 1663
$ws.Range("K113").Value = 1734.8
$ws.Range("L113").Value = 1663
$ws.Range("M113").Value = 435.2
$ws.Range("N113").Value = -6003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 113295
$ws.Range("I132").Value = 1275.3334
$ws.Range("J132").Value = 337334.34
$ws.Range("K132").Value = 3826.0002
$ws.Range("L132").Value = 1012003.02
$ws.Range("M132").Value = -1296.0002
$ws.Range("N132").Value = -1017063.02

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 37300
$ws.Range("J133").Value = 37300
$ws.Range("L133").Value = 37300
$ws.Range("N133").Value = -42360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1998
$ws.Range("I136").Value = 1947.3334
$ws.Range("J136").Value = 2150
$ws.Range("K136").Value = 5842.0002
$ws.Range("L136").Value = 6450
$ws.Range("M136").Value = -3292.0002
$ws.Range("N136").Value = -11550

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 739.5454999999999
$ws.Range("I113").Value = 355.33334
$ws.Range("J113").Value = 1200.6
$ws.Range("K113").Value = 1066.00002
$ws.Range("L113").Value = 3601.8
$ws.Range("M113").Value = 1103.99998
$ws.Range("N113").Value = -7941.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4416.706
$ws.Range("I132").Value = 4174.25
$ws.Range("K132").Value = 12522.75
$ws.Range("M132").Value = -9992.75
